$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Worksheet 1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 194
$ws1.Range("F5").Value = 5188
$ws1.Range("F7").Value = 50
$ws1.Range("F9").Value = 580
$ws1.Range("F10").Value = 533
$ws1.Range("F13").Value = 1433
$ws1.Range("F14").Value = 4135
$ws1.Range("F15").Value = 426
$ws1.Range("F16").Value = 167
$ws1.Range("F17").Value = 146
$ws1.Range("F18").Value = 92
$ws1.Range("F19").Value = 3141
$ws1.Range("F21").Value = 1058
$ws1.Range("F24").Value = 187
$ws1.Range("F25").Value = 95
$ws1.Range("F26").Value = 25
$ws1.Range("F29").Value = 286
$ws1.Range("F30").Value = 17
$ws1.Range("F33").Value = 11
$ws1.Range("F34").Value = 5

# --- Sheet "全部类型" (Worksheet 4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 194
$ws4.Range("F6").Value = 5188
$ws4.Range("F8").Value = 50
$ws4.Range("F10").Value = 580
$ws4.Range("F11").Value = 533
$ws4.Range("F14").Value = 1433
$ws4.Range("F15").Value = 4136
$ws4.Range("F16").Value = 426
$ws4.Range("F17").Value = 167
$ws4.Range("F18").Value = 146
$ws4.Range("F19").Value = 92
$ws4.Range("F20").Value = 3142
$ws4.Range("F22").Value = 1058
$ws4.Range("F25").Value = 187
$ws4.Range("F26").Value = 95
$ws4.Range("F27").Value = 25
$ws4.Range("F30").Value = 286
$ws4.Range("F31").Value = 17
$ws4.Range("F34").Value = 11
$ws4.Range("F35").Value = 5

$wb.Save()
